# edit.ps1
#
# Daily "BP terminal gate pricing by state" data refresh.
#
# The workbook lists, for each terminal, the two most recent effective-date
# pricing rows (newest on top). This edit performs the routine daily rollover:
#   - the pricing row that was previously "today" (effective date serial
#     45951 = 21-Oct-2025) slides down into the "previous day" slot, replacing
#     the now-stale 18-Oct-2025 (serial 45948) entry in that block, and
#   - the top row of each terminal block is populated with the newly
#     published 22-Oct-2025 (serial 45952) pricing.
# Only the Effective Date (col A) and the Diesel/ULP/PULP/e10 price columns
# (D/E/F/G) change; Terminal names/styles/columns B & C are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 8: Sydney-Botany
    $ws.Cells.Item(8, 1).Value = 45952
    $ws.Cells.Item(8, 4).Value = 160.06
    $ws.Cells.Item(8, 5).Value = 157.88999999999999
    $ws.Cells.Item(8, 6).Value = 167.89
    $ws.Cells.Item(8, 7).Value = 158.05000000000001
    # Row 9: Sydney-Silverwater
    $ws.Cells.Item(9, 1).Value = 45952
    $ws.Cells.Item(9, 4).Value = 160.06
    $ws.Cells.Item(9, 5).Value = 157.88999999999999
    $ws.Cells.Item(9, 6).Value = 167.89
    $ws.Cells.Item(9, 7).Value = 158.05000000000001
    # Row 10: Newcastle
    $ws.Cells.Item(10, 1).Value = 45952
    $ws.Cells.Item(10, 4).Value = 162.28
    $ws.Cells.Item(10, 5).Value = 160.22999999999999
    $ws.Cells.Item(10, 6).Value = 170.23
    $ws.Cells.Item(10, 7).Value = 160.69999999999999
    # Row 11: Sydney-Botany
    $ws.Cells.Item(11, 1).Value = 45951
    $ws.Cells.Item(11, 4).Value = 160.47
    $ws.Cells.Item(11, 5).Value = 157.9
    $ws.Cells.Item(11, 6).Value = 167.9
    $ws.Cells.Item(11, 7).Value = 158.06
    # Row 12: Sydney-Silverwater
    $ws.Cells.Item(12, 1).Value = 45951
    $ws.Cells.Item(12, 4).Value = 160.47
    $ws.Cells.Item(12, 5).Value = 157.9
    $ws.Cells.Item(12, 6).Value = 167.9
    $ws.Cells.Item(12, 7).Value = 158.06
    # Row 13: Newcastle
    $ws.Cells.Item(13, 1).Value = 45951
    $ws.Cells.Item(13, 4).Value = 162.69
    $ws.Cells.Item(13, 5).Value = 160.26
    $ws.Cells.Item(13, 6).Value = 170.26
    $ws.Cells.Item(13, 7).Value = 160.72999999999999
    # Row 17: Darwin
    $ws.Cells.Item(17, 1).Value = 45952
    $ws.Cells.Item(17, 4).Value = 165.7
    $ws.Cells.Item(17, 5).Value = 163.07
    $ws.Cells.Item(17, 6).Value = 173.07
    # Row 18: Darwin
    $ws.Cells.Item(18, 1).Value = 45951
    $ws.Cells.Item(18, 4).Value = 166.12
    $ws.Cells.Item(18, 5).Value = 163.09
    $ws.Cells.Item(18, 6).Value = 173.09
    # Row 22: Brisbane
    $ws.Cells.Item(22, 1).Value = 45952
    $ws.Cells.Item(22, 4).Value = 160.97999999999999
    $ws.Cells.Item(22, 5).Value = 159.12
    $ws.Cells.Item(22, 6).Value = 168.72
    $ws.Cells.Item(22, 7).Value = 160.30000000000001
    # Row 23: Cairns
    $ws.Cells.Item(23, 1).Value = 45952
    $ws.Cells.Item(23, 4).Value = 167.04
    $ws.Cells.Item(23, 5).Value = 163.83000000000001
    $ws.Cells.Item(23, 6).Value = 173.83
    # Row 24: Gladstone
    $ws.Cells.Item(24, 1).Value = 45952
    $ws.Cells.Item(24, 4).Value = 166.85
    $ws.Cells.Item(24, 5).Value = 164.04
    $ws.Cells.Item(24, 6).Value = 174.04
    # Row 25: Mackay
    $ws.Cells.Item(25, 1).Value = 45952
    $ws.Cells.Item(25, 4).Value = 167.68
    $ws.Cells.Item(25, 5).Value = 163.43
    $ws.Cells.Item(25, 6).Value = 173.43
    $ws.Cells.Item(25, 7).Value = 163.26
    # Row 26: Townsville
    $ws.Cells.Item(26, 1).Value = 45952
    $ws.Cells.Item(26, 4).Value = 166.4
    $ws.Cells.Item(26, 5).Value = 164.97
    $ws.Cells.Item(26, 6).Value = 174.97
    # Row 27: Brisbane
    $ws.Cells.Item(27, 1).Value = 45951
    $ws.Cells.Item(27, 4).Value = 161.38999999999999
    $ws.Cells.Item(27, 5).Value = 159.15
    $ws.Cells.Item(27, 6).Value = 168.75
    $ws.Cells.Item(27, 7).Value = 160.32
    # Row 28: Cairns
    $ws.Cells.Item(28, 1).Value = 45951
    $ws.Cells.Item(28, 4).Value = 167.46
    $ws.Cells.Item(28, 5).Value = 163.86
    $ws.Cells.Item(28, 6).Value = 173.86
    # Row 29: Gladstone
    $ws.Cells.Item(29, 1).Value = 45951
    $ws.Cells.Item(29, 4).Value = 167.26
    $ws.Cells.Item(29, 5).Value = 164.07
    $ws.Cells.Item(29, 6).Value = 174.07
    # Row 30: Mackay
    $ws.Cells.Item(30, 1).Value = 45951
    $ws.Cells.Item(30, 4).Value = 168.1
    $ws.Cells.Item(30, 5).Value = 163.46
    $ws.Cells.Item(30, 6).Value = 173.46
    $ws.Cells.Item(30, 7).Value = 163.29
    # Row 31: Townsville
    $ws.Cells.Item(31, 1).Value = 45951
    $ws.Cells.Item(31, 4).Value = 166.81
    $ws.Cells.Item(31, 5).Value = 165
    $ws.Cells.Item(31, 6).Value = 175
    # Row 35: Adelaide
    $ws.Cells.Item(35, 1).Value = 45952
    $ws.Cells.Item(35, 4).Value = 160.52000000000001
    $ws.Cells.Item(35, 5).Value = 157.33000000000001
    $ws.Cells.Item(35, 6).Value = 166.33
    # Row 36: Adelaide
    $ws.Cells.Item(36, 1).Value = 45951
    $ws.Cells.Item(36, 4).Value = 161.05000000000001
    $ws.Cells.Item(36, 5).Value = 157.36000000000001
    $ws.Cells.Item(36, 6).Value = 166.36
    # Row 40: Burnie
    $ws.Cells.Item(40, 1).Value = 45952
    $ws.Cells.Item(40, 4).Value = 166.17
    $ws.Cells.Item(40, 5).Value = 162.81
    $ws.Cells.Item(40, 6).Value = 172.81
    # Row 41: Hobart
    $ws.Cells.Item(41, 1).Value = 45952
    $ws.Cells.Item(41, 4).Value = 165.89
    $ws.Cells.Item(41, 5).Value = 163.22999999999999
    $ws.Cells.Item(41, 6).Value = 173.23
    # Row 42: Burnie
    $ws.Cells.Item(42, 1).Value = 45951
    $ws.Cells.Item(42, 4).Value = 166.6
    $ws.Cells.Item(42, 5).Value = 162.84
    $ws.Cells.Item(42, 6).Value = 172.84
    # Row 43: Hobart
    $ws.Cells.Item(43, 1).Value = 45951
    $ws.Cells.Item(43, 4).Value = 166.31
    $ws.Cells.Item(43, 5).Value = 163.26
    $ws.Cells.Item(43, 6).Value = 173.26
    # Row 47: Geelong
    $ws.Cells.Item(47, 1).Value = 45952
    $ws.Cells.Item(47, 4).Value = 160.86000000000001
    $ws.Cells.Item(47, 5).Value = 159.02000000000001
    $ws.Cells.Item(47, 6).Value = 169.02
    # Row 48: Melbourne
    $ws.Cells.Item(48, 1).Value = 45952
    $ws.Cells.Item(48, 4).Value = 160.84
    $ws.Cells.Item(48, 5).Value = 159.19
    $ws.Cells.Item(48, 6).Value = 169.19
    # Row 49: Geelong
    $ws.Cells.Item(49, 1).Value = 45951
    $ws.Cells.Item(49, 4).Value = 161.34
    $ws.Cells.Item(49, 5).Value = 159.06
    $ws.Cells.Item(49, 6).Value = 169.06
    # Row 50: Melbourne
    $ws.Cells.Item(50, 1).Value = 45951
    $ws.Cells.Item(50, 4).Value = 161.32
    $ws.Cells.Item(50, 5).Value = 159.22999999999999
    $ws.Cells.Item(50, 6).Value = 169.23
    # Row 54: Broome
    $ws.Cells.Item(54, 1).Value = 45952
    $ws.Cells.Item(54, 4).Value = 176.36
    $ws.Cells.Item(54, 5).Value = 173.11
    $ws.Cells.Item(54, 6).Value = 183.11
    # Row 55: Esperance
    $ws.Cells.Item(55, 1).Value = 45952
    $ws.Cells.Item(55, 4).Value = 164.01
    $ws.Cells.Item(55, 5).Value = 170.5
    $ws.Cells.Item(55, 6).Value = 180.5
    # Row 56: Geraldton
    $ws.Cells.Item(56, 1).Value = 45952
    $ws.Cells.Item(56, 4).Value = 166.3
    # Row 57: Kalgoorlie
    $ws.Cells.Item(57, 1).Value = 45952
    $ws.Cells.Item(57, 4).Value = 165.96
    $ws.Cells.Item(57, 5).Value = 164.77
    # Row 58: Perth
    $ws.Cells.Item(58, 1).Value = 45952
    $ws.Cells.Item(58, 4).Value = 161.87
    $ws.Cells.Item(58, 5).Value = 160.82
    $ws.Cells.Item(58, 6).Value = 170.82
    # Row 59: Port Hedland
    $ws.Cells.Item(59, 1).Value = 45952
    $ws.Cells.Item(59, 4).Value = 168.69
    $ws.Cells.Item(59, 5).Value = 171.29
    # Row 60: Broome
    $ws.Cells.Item(60, 1).Value = 45951
    $ws.Cells.Item(60, 4).Value = 176.77
    $ws.Cells.Item(60, 5).Value = 173.16
    $ws.Cells.Item(60, 6).Value = 183.16
    # Row 61: Esperance
    $ws.Cells.Item(61, 1).Value = 45951
    $ws.Cells.Item(61, 4).Value = 164.43
    $ws.Cells.Item(61, 5).Value = 170.52
    $ws.Cells.Item(61, 6).Value = 180.52
    # Row 62: Geraldton
    $ws.Cells.Item(62, 1).Value = 45951
    $ws.Cells.Item(62, 4).Value = 166.82
    # Row 63: Kalgoorlie
    $ws.Cells.Item(63, 1).Value = 45951
    $ws.Cells.Item(63, 4).Value = 166.48
    $ws.Cells.Item(63, 5).Value = 164.79
    # Row 64: Perth
    $ws.Cells.Item(64, 1).Value = 45951
    $ws.Cells.Item(64, 4).Value = 162.38999999999999
    $ws.Cells.Item(64, 5).Value = 160.84
    $ws.Cells.Item(64, 6).Value = 170.84
    # Row 65: Port Hedland
    $ws.Cells.Item(65, 1).Value = 45951
    $ws.Cells.Item(65, 4).Value = 169.1
    $ws.Cells.Item(65, 5).Value = 171.33
